$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.617.36"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "3.323.57"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "579.72"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "174.63"
$ws.Range("E6").Value = "  -4.40%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "3.320.54"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "0.575"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "45.23"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").Value = "662.06"
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("D15").Value = "3.866.89"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "67.628.45"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "3.326.58"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "17.32"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "5.36"
$ws.Range("E23").Value = "  +6.82%  "
$ws.Range("D24").Value = "17.02"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").Value = "98.75"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "9.22"
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("D29").Value = "33.38"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").Value = "7.28"
$ws.Range("E31").Value = "  +9.30%  "
$ws.Range("D32").Value = "570.28"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.673.21"
$ws.Range("E36").Value = "  -6.74%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "56.54"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  -6.76%  "
$ws.Range("D39").Value = "34.14"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "3.35"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "3.10"
$ws.Range("E43").Value = "  -4.58%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Value = "0.0₃0661"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "129.58"
$ws.Range("E51").Value = "  -1.04%  "
